$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.818807045481801
$ws.Range("D2").Value = 9.161649988426518
$ws.Range("E2").Value = 13.67796320441161
$ws.Range("F2").Value = 33.13432706867357
$ws.Range("G2").Value = 3.657624966194494
$ws.Range("J2").Value = 9.943717855007101
$ws.Range("K2").Value = 12.87557468250796
$ws.Range("O2").Value = 25.00234414343934

$ws.Range("B3").Value = 7.748103423656503
$ws.Range("D3").Value = 9.093153101097551
$ws.Range("E3").Value = 13.61111433209262
$ws.Range("F3").Value = 33.21631978249429
$ws.Range("G3").Value = 3.659817400840919
$ws.Range("J3").Value = 9.951689377248785
$ws.Range("K3").Value = 12.41343131606949
$ws.Range("O3").Value = 25.11078579430545

$ws.Range("B4").Value = 7.7062190359354
$ws.Range("D4").Value = 9.05219716643666
$ws.Range("E4").Value = 13.57268373807857
$ws.Range("F4").Value = 33.27606293987089
$ws.Range("G4").Value = 3.661234502412696
$ws.Range("J4").Value = 9.958219330258199
$ws.Range("K4").Value = 12.12120271468018
$ws.Range("O4").Value = 25.18384123941853

$ws.Range("B5").Value = 7.689553792019252
$ws.Range("D5").Value = 9.03579801745045
$ws.Range("E5").Value = 13.55769263197289
$ws.Range("F5").Value = 33.30276389229012
$ws.Range("G5").Value = 3.661829878419572
$ws.Range("J5").Value = 9.961291725580431
$ws.Range("K5").Value = 12.00015017795923
$ws.Range("O5").Value = 25.21523403794422

$ws.Range("B6").Value = 7.68681141416403
$ws.Range("D6").Value = 9.033092903206466
$ws.Range("E6").Value = 13.55524415362335
$ws.Range("F6").Value = 33.30733954942804
$ws.Range("G6").Value = 3.661929822700615
$ws.Range("J6").Value = 9.961826743501025
$ws.Range("K6").Value = 11.97993578715329
$ws.Range("O6").Value = 25.22054460393018

$ws.Range("B7").Value = 7.705992626704377
$ws.Range("D7").Value = 9.051974806803125
$ws.Range("E7").Value = 13.57247883593335
$ws.Range("F7").Value = 33.27641351410345
$ws.Range("G7").Value = 3.661242459326931
$ws.Range("J7").Value = 9.958259099849782
$ws.Range("K7").Value = 12.11957789447645
$ws.Range("O7").Value = 25.18425805288345

$ws.Range("B8").Value = 7.794122802631111
$ws.Range("D8").Value = 9.137812056055832
$ws.Range("E8").Value = 13.65437854804483
$ws.Range("F8").Value = 33.1606424973141
$ws.Range("G8").Value = 3.658366229320456
$ws.Range("J8").Value = 9.94612712693743
$ws.Range("K8").Value = 12.71807840214115
$ws.Range("O8").Value = 25.03838844256347

$ws.Range("B9").Value = 7.978210935358105
$ws.Range("D9").Value = 9.314253306572825
$ws.Range("E9").Value = 13.8351577585505
$ws.Range("F9").Value = 33.00853926349052
$ws.Range("G9").Value = 3.653286174774896
$ws.Range("J9").Value = 9.935302737688733
$ws.Range("K9").Value = 13.81819641536491
$ws.Range("O9").Value = 24.80392281191325

$ws.Range("B10").Value = 8.119140351964898
$ws.Range("D10").Value = 9.447970796441515
$ws.Range("E10").Value = 13.97945161241797
$ws.Range("F10").Value = 32.94288717778409
$ws.Range("G10").Value = 3.649891692488063
$ws.Range("J10").Value = 9.935236608243585
$ws.Range("K10").Value = 14.57424931413739
$ws.Range("O10").Value = 24.66342218051234

$ws.Range("B11").Value = 8.184226159695696
$ws.Range("D11").Value = 9.509500860029702
$ws.Range("E11").Value = 14.04739554121851
$ws.Range("F11").Value = 32.92310303001462
$ws.Range("G11").Value = 3.648420025645579
$ws.Range("J11").Value = 9.93691284073528
$ws.Range("K11").Value = 14.90559438097017
$ws.Range("O11").Value = 24.6064669924057

$ws.Range("B12").Value = 8.20899084244556
$ws.Range("D12").Value = 9.532885407654955
$ws.Range("E12").Value = 14.07343836433624
$ws.Range("F12").Value = 32.91706540453803
$ws.Range("G12").Value = 3.647873109499343
$ws.Range("J12").Value = 9.937792194919904
$ws.Range("K12").Value = 15.02916780033541
$ws.Range("O12").Value = 24.58590555394339

$ws.Range("B13").Value = 8.203652426900499
$ws.Range("D13").Value = 9.527845632732403
$ws.Range("E13").Value = 14.0678158964739
$ws.Range("F13").Value = 32.91830097763614
$ws.Range("G13").Value = 3.647990437217215
$ws.Range("J13").Value = 9.937591945078205
$ws.Range("K13").Value = 15.00263985146572
$ws.Range("O13").Value = 24.59028898756594

$ws.Range("B14").Value = 8.186261325705754
$ws.Range("D14").Value = 9.511423122494268
$ws.Range("E14").Value = 14.04953190969069
$ws.Range("F14").Value = 32.92257714415101
$ws.Range("G14").Value = 3.648374822955568
$ws.Range("J14").Value = 9.936980288388634
$ws.Range("K14").Value = 14.91579938730411
$ws.Range("O14").Value = 24.60475519246168

$ws.Range("B15").Value = 8.175623493102472
$ws.Range("D15").Value = 9.501374370733974
$ws.Range("E15").Value = 14.03837279828726
$ws.Range("F15").Value = 32.92538591709013
$ws.Range("G15").Value = 3.648611619504795
$ws.Range("J15").Value = 9.936637460460485
$ws.Range("K15").Value = 14.86235722613666
$ws.Range("O15").Value = 24.61374736650589

$ws.Range("B16").Value = 8.114904581291679
$ws.Range("D16").Value = 9.443962373851553
$ws.Range("E16").Value = 13.97505620558913
$ws.Range("F16").Value = 32.9443833889426
$ws.Range("G16").Value = 3.649989323702086
$ws.Range("J16").Value = 9.935161335450786
$ws.Range("K16").Value = 14.55233390649533
$ws.Range("O16").Value = 24.66728484899336

$ws.Range("B17").Value = 8.077889691334443
$ws.Range("D17").Value = 9.408909950876227
$ws.Range("E17").Value = 13.93679143807509
$ws.Range("F17").Value = 32.958623314603
$ws.Range("G17").Value = 3.650853031596359
$ws.Range("J17").Value = 9.934692304864127
$ws.Range("K17").Value = 14.35885311229178
$ws.Range("O17").Value = 24.70191476180958

$ws.Range("B18").Value = 8.056692937941252
$ws.Range("D18").Value = 9.388815722611561
$ws.Range("E18").Value = 13.91500049010073
$ws.Range("F18").Value = 32.9677624346926
$ws.Range("G18").Value = 3.651356640753114
$ws.Range("J18").Value = 9.934583191084812
$ws.Range("K18").Value = 14.24638788333968
$ws.Range("O18").Value = 24.72248768688354

$ws.Range("B19").Value = 8.049532773517708
$ws.Range("D19").Value = 9.382024167710682
$ws.Range("E19").Value = 13.9076604000831
$ws.Range("F19").Value = 32.97101956876088
$ws.Range("G19").Value = 3.651528328372071
$ws.Range("J19").Value = 9.934573861809943
$ws.Range("K19").Value = 14.20810931927898
$ws.Range("O19").Value = 24.72956560348209

$ws.Range("B20").Value = 8.081820501909343
$ws.Range("D20").Value = 9.412634525660877
$ws.Range("E20").Value = 13.94084235661794
$ws.Range("F20").Value = 32.95700922830937
$ws.Range("G20").Value = 3.650760382212
$ws.Range("J20").Value = 9.934725611379365
$ws.Range("K20").Value = 14.37957230616673
$ws.Range("O20").Value = 24.6981605474628

$ws.Range("B21").Value = 8.191366485968329
$ws.Range("D21").Value = 9.516244650902081
$ws.Range("E21").Value = 14.05489398089435
$ws.Range("F21").Value = 32.92128163259565
$ws.Range("G21").Value = 3.648261638472221
$ws.Range("J21").Value = 9.937153315098172
$ws.Range("K21").Value = 14.94135875319445
$ws.Range("O21").Value = 24.60047876456416

$ws.Range("B22").Value = 8.263639207896846
$ws.Range("D22").Value = 9.584444455376564
$ws.Range("E22").Value = 14.13125472698086
$ws.Range("F22").Value = 32.90640918451861
$ws.Range("G22").Value = 3.646688997355205
$ws.Range("J22").Value = 9.940165239204783
$ws.Range("K22").Value = 15.29741175340821
$ws.Range("O22").Value = 24.54250576676509

$ws.Range("B23").Value = 8.225011199992997
$ws.Range("D23").Value = 9.548006023406517
$ws.Range("E23").Value = 14.09033877325355
$ws.Range("F23").Value = 32.91356994993859
$ws.Range("G23").Value = 3.647522833456698
$ws.Range("J23").Value = 9.938427590858772
$ws.Range("K23").Value = 15.10842253496534
$ws.Range("O23").Value = 24.57290836161072

$ws.Range("B24").Value = 8.080043120354372
$ws.Range("D24").Value = 9.410950463889019
$ws.Range("E24").Value = 13.93901028722384
$ws.Range("F24").Value = 32.95773599037054
$ws.Range("G24").Value = 3.65080224702793
$ws.Range("J24").Value = 9.934710053386306
$ws.Range("K24").Value = 14.37020898455151
$ws.Range("O24").Value = 24.69985576066521

$ws.Range("B25").Value = 7.927325191428943
$ws.Range("D25").Value = 9.265740199638186
$ws.Range("E25").Value = 13.78417431109295
$ws.Range("F25").Value = 33.04161914335989
$ws.Range("G25").Value = 3.654600871587271
$ws.Range("J25").Value = 9.936844280135233
$ws.Range("K25").Value = 13.52927743467873
$ws.Range("O25").Value = 24.86179673434239
